# Update "想去人数" (number of people interested) figures for several
# events on the "展览" (Exhibition) and "全部类型" (All Types) sheets.
# Source: gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions only) ---
$wsExhibit = $wb.Sheets("展览")
$wsExhibit.Range("F2").Value  = 4517
$wsExhibit.Range("F4").Value  = 141
$wsExhibit.Range("F8").Value  = 636
$wsExhibit.Range("F10").Value = 196
$wsExhibit.Range("F11").Value = 1376
$wsExhibit.Range("F12").Value = 28
$wsExhibit.Range("F13").Value = 2976
$wsExhibit.Range("F14").Value = 448
$wsExhibit.Range("F15").Value = 671

# --- Sheet "全部类型" (combined listing, rows offset by +1 starting row 10
#     because it also includes the "本地生活" performance entry) ---
$wsAll = $wb.Sheets("全部类型")
$wsAll.Range("F2").Value  = 4517
$wsAll.Range("F4").Value  = 141
$wsAll.Range("F8").Value  = 636
$wsAll.Range("F11").Value = 196
$wsAll.Range("F12").Value = 1376
$wsAll.Range("F13").Value = 28
$wsAll.Range("F14").Value = 2976
$wsAll.Range("F15").Value = 448
$wsAll.Range("F16").Value = 671
